$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Soldatino di Pb"
$ws.Range("B3").Value = "Stefano Pizzini"
$ws.Range("C3").Value = "Matteo Mazzola"
$ws.Range("D3").Value = "Davide Raffaelli"
$ws.Range("E3").Value = "Michele Parisi"
$ws.Range("F3").Value = "Leonardo Parisi"
